$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.963.27"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "2.912.95"
$ws.Range("E3").Value = "  -4.03%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "586.28"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").Value = "146.53"
$ws.Range("E6").Value = "  -3.53%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -3.42%  "
$ws.Range("D9").Value = "2.910.39"
$ws.Range("E9").Value = "  -3.97%  "
$ws.Range("D10").Value = "6.78"
$ws.Range("E10").Value = "  +6.74%  "
$ws.Range("E11").Value = "  -4.62%  "
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("D13").Value = "0.0000225"
$ws.Range("E13").Value = "  -4.08%  "
$ws.Range("D14").Value = "33.64"
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "3.394.09"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").Value = "60.885.64"
$ws.Range("E17").Value = "  -3.20%  "
$ws.Range("D18").Value = "6.76"
$ws.Range("E18").Value = "  -3.83%  "
$ws.Range("D19").Value = "2.909.50"
$ws.Range("E19").Value = "  -4.18%  "
$ws.Range("D20").Value = "427.24"
$ws.Range("E20").Value = "  -6.29%  "
$ws.Range("D21").Value = "13.60"
$ws.Range("E21").Value = "  -4.34%  "
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("D23").Value = "7.12"
$ws.Range("E23").Value = "  -4.94%  "
$ws.Range("D24").Value = "80.29"
$ws.Range("E24").Value = "  -3.30%  "
$ws.Range("D25").Value = "10.98"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -2.25%  "
$ws.Range("D27").Value = "11.88"
$ws.Range("E27").Value = "  -2.15%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "7.25"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "26.47"
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("E34").Value = "  -4.11%  "
$ws.Range("D35").Value = "0.0₃0846"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("D37").Value = "5.63"
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("D38").Value = "2.99"
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D39").Value = "49.34"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("D42").Value = "8.67"
$ws.Range("E42").Value = "  -4.87%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "41.61"
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("D45").Value = "377.14"
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("D46").Value = "0.0347"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("D47").Value = "2.671.37"
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").Value = "133.06"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D50").Value = "24.78"
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("E51").Value = "  -1.56%  "
